# Update example input files to include Dp
#
# Adds a new "Dp" / "Dp_units" pair of columns to the "ions" sheet, inserted
# between the existing "Ds_units" column (H) and the trailing "conc_units"
# column (which shifts from I to K).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ions")

# Insert two new blank columns at I:J. This pushes the old "conc_units"
# column (and its data) from column I to column K.
$ws.Range("I1:J1").EntireColumn.Insert()

# Header row
$ws.Range("I1").Value = "Dp"
$ws.Range("J1").Value = "Dp_units"

# Dp data: row 2 (CHLORIDE) = 1, rows 3-6 = 2e-6 shown in scientific notation
# (matching the display style already used by the kL / Ds columns).
$ws.Range("I2").Value = 1
$ws.Range("I3:I6").Value = 0.000002
$ws.Range("I3:I6").NumberFormat = "0.00E+00"

# Dp_units data: cm^2/s for every row, same unit text as Ds_units.
$ws.Range("J2:J6").Value = "cm^2/s"

# The "ions" sheet becomes the active tab (previously "params" was active).
$ws.Activate()
$null = $ws.Range("A1").Select()
